# Auto-generated Excel COM-interop script to apply market-data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1001154.9
$ws.Range("I2").Value = 683.3333
$ws.Range("J2").Value = 1429928.4
$ws.Range("K2").Value = 683.3333
$ws.Range("L2").Value = 1429928.4
$ws.Range("M2").Value = -570.3333
$ws.Range("N2").Value = -1430154.4
$ws.Range("H4").Value = 1158.75
$ws.Range("I4").Value = 805.75
$ws.Range("J4").Value = 1688.25
$ws.Range("K4").Value = 805.75
$ws.Range("L4").Value = 1688.25
$ws.Range("M4").Value = -691.75
$ws.Range("N4").Value = -1916.25
$ws.Range("H40").Value = 55558320
$ws.Range("I40").Value = 3050
$ws.Range("K40").Value = 3050
$ws.Range("M40").Value = -2875
$ws.Range("H135").Value = 2231.7778
$ws.Range("I135").Value = 1075.5
$ws.Range("J135").Value = 4544.3335
$ws.Range("K135").Value = 9679.5
$ws.Range("L135").Value = 40899.0015
$ws.Range("M135").Value = -7144.5
$ws.Range("N135").Value = -45969.0015
$ws.Range("H138").Value = 2778.2917
$ws.Range("I138").Value = 2861.762
$ws.Range("J138").Value = 2194
$ws.Range("K138").Value = 8585.286
$ws.Range("L138").Value = 6582
$ws.Range("M138").Value = -3445.286
$ws.Range("N138").Value = -16862

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2144.647
$ws.Range("I2").Value = 2244.3704
$ws.Range("J2").Value = 1760
$ws.Range("K2").Value = 2244.3704
$ws.Range("L2").Value = 1760
$ws.Range("M2").Value = -2131.3704
$ws.Range("N2").Value = -1986
$ws.Range("H45").Value = 3412.5334
$ws.Range("I45").Value = 2622.111
$ws.Range("J45").Value = 4598.1665
$ws.Range("K45").Value = 2622.111
$ws.Range("L45").Value = 4598.1665
$ws.Range("M45").Value = -2245.111
$ws.Range("N45").Value = -5352.1665
$ws.Range("H88").Value = 1989.8182
$ws.Range("I88").Value = 1231.2
$ws.Range("J88").Value = 2622
$ws.Range("K88").Value = 1231.2
$ws.Range("L88").Value = 2622
$ws.Range("M88").Value = -825.2
$ws.Range("N88").Value = -3434
$ws.Range("H91").Value = 1989.8182
$ws.Range("I91").Value = 1231.2
$ws.Range("J91").Value = 2622
$ws.Range("K91").Value = 1231.2
$ws.Range("L91").Value = 2622
$ws.Range("M91").Value = 172.8
$ws.Range("N91").Value = -5430
$ws.Range("H116").Value = 2144.647
$ws.Range("I116").Value = 2244.3704
$ws.Range("J116").Value = 1760
$ws.Range("K116").Value = 2244.3704
$ws.Range("L116").Value = 1760
$ws.Range("M116").Value = 49.62960000000021
$ws.Range("N116").Value = -6348

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2144.647
$ws.Range("I3").Value = 2244.3704
$ws.Range("J3").Value = 1760
$ws.Range("K3").Value = 2244.3704
$ws.Range("L3").Value = 1760
$ws.Range("M3").Value = -2130.3704
$ws.Range("N3").Value = -1988
$ws.Range("H99").Value = 1555.44
$ws.Range("I99").Value = 1566.0476
$ws.Range("K99").Value = 1566.0476
$ws.Range("M99").Value = -68.0476000000001
$ws.Range("H134").Value = 3207.6924
$ws.Range("I134").Value = 2265.8
$ws.Range("J134").Value = 6347.3335
$ws.Range("K134").Value = 6797.400000000001
$ws.Range("L134").Value = 19042.0005
$ws.Range("M134").Value = -4262.400000000001
$ws.Range("N134").Value = -24112.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38464188
$ws.Range("I31").Value = 58825908
$ws.Range("J31").Value = 3166.2222
$ws.Range("K31").Value = 58825908
$ws.Range("L31").Value = 3166.2222
$ws.Range("M31").Value = -58825613
$ws.Range("N31").Value = -3756.2222
$ws.Range("H34").Value = 38464188
$ws.Range("I34").Value = 58825908
$ws.Range("J34").Value = 3166.2222
$ws.Range("K34").Value = 58825908
$ws.Range("L34").Value = 3166.2222
$ws.Range("M34").Value = -58825706
$ws.Range("N34").Value = -3570.2222
$ws.Range("H99").Value = 23447.4
$ws.Range("J99").Value = 23830.143
$ws.Range("L99").Value = 23830.143
$ws.Range("N99").Value = -26826.143
$ws.Range("H107").Value = 996.2
$ws.Range("I107").Value = 687.3333
$ws.Range("J107").Value = 2617.75
$ws.Range("K107").Value = 687.3333
$ws.Range("L107").Value = 2617.75
$ws.Range("M107").Value = 1232.6667
$ws.Range("N107").Value = -6457.75
$ws.Range("H122").Value = 3515.3635
$ws.Range("I122").Value = 3189.6428
$ws.Range("J122").Value = 4085.375
$ws.Range("K122").Value = 9568.928400000001
$ws.Range("L122").Value = 12256.125
$ws.Range("M122").Value = -7118.928400000001
$ws.Range("N122").Value = -17156.125
$ws.Range("H126").Value = 23447.4
$ws.Range("J126").Value = 23830.143
$ws.Range("L126").Value = 71490.429
$ws.Range("N126").Value = -76430.429
$ws.Range("H132").Value = 2788.88
$ws.Range("I132").Value = 2506.4375
$ws.Range("J132").Value = 3291
$ws.Range("K132").Value = 7519.3125
$ws.Range("L132").Value = 9873
$ws.Range("M132").Value = -4989.3125
$ws.Range("N132").Value = -14933
$ws.Range("H134").Value = 1885.8572
$ws.Range("I134").Value = 1920.4
$ws.Range("K134").Value = 5761.200000000001
$ws.Range("M134").Value = -3226.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2337.4167
$ws.Range("I11").Value = 424
$ws.Range("K11").Value = 1272
$ws.Range("M11").Value = -1132
$ws.Range("H69").Value = 9186.4
$ws.Range("I69").Value = 3199.6667
$ws.Range("K69").Value = 9599.000100000001
$ws.Range("M69").Value = -8788.000100000001
$ws.Range("H72").Value = 9186.4
$ws.Range("I72").Value = 3199.6667
$ws.Range("K72").Value = 28797.0003
$ws.Range("M72").Value = -24741.0003
$ws.Range("H98").Value = 718.75
$ws.Range("I98").Value = 725
$ws.Range("J98").Value = 712.5
$ws.Range("K98").Value = 2175
$ws.Range("L98").Value = 2137.5
$ws.Range("M98").Value = -677
$ws.Range("N98").Value = -5133.5
$ws.Range("H107").Value = 4337225
$ws.Range("J107").Value = 5691651.5
$ws.Range("L107").Value = 17074954.5
$ws.Range("N107").Value = -17078794.5
$ws.Range("H132").Value = 2498
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 39108.57
$ws.Range("I133").Value = 40071.168
$ws.Range("J133").Value = 33333
$ws.Range("K133").Value = 120213.504
$ws.Range("L133").Value = 99999
$ws.Range("M133").Value = -115153.504
$ws.Range("N133").Value = -110119

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2333
$ws.Range("I97").Value = 2999
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 2999
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -2503
$ws.Range("N97").Value = -2992
$ws.Range("H122").Value = 2128.111
$ws.Range("I122").Value = 1066.5
$ws.Range("J122").Value = 2977.4
$ws.Range("K122").Value = 3199.5
$ws.Range("L122").Value = 8932.200000000001
$ws.Range("M122").Value = -749.5
$ws.Range("N122").Value = -13832.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3099.348
$ws.Range("I40").Value = 2814.75
$ws.Range("K40").Value = 2814.75
$ws.Range("M40").Value = -2678.75
$ws.Range("H46").Value = 2699.4
$ws.Range("J46").Value = 2999.25
$ws.Range("L46").Value = 2999.25
$ws.Range("N46").Value = -3375.25
$ws.Range("H132").Value = 2040.8948
$ws.Range("J132").Value = 2334.6667
$ws.Range("L132").Value = 7004.000100000001
$ws.Range("N132").Value = -12064.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4068.361
$ws.Range("I107").Value = 2147.0667
$ws.Range("K107").Value = 6441.2001
$ws.Range("M107").Value = -4521.2001
